$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1256.8276
$ws.Range("J17").Value = 1256.8276
$ws.Range("L17").Value = 3770.4828
$ws.Range("N17").Value = -4106.4828

$ws.Range("H29").Value = 3335.5715
$ws.Range("J29").Value = 5649.75
$ws.Range("L29").Value = 16949.25
$ws.Range("N29").Value = -17511.25

$ws.Range("H43").Value = 2494.5
$ws.Range("J43").Value = 2494.5
$ws.Range("L43").Value = 2494.5
$ws.Range("N43").Value = -2632.5

$ws.Range("H51").Value = 17221.857
$ws.Range("I51").Value = 37500
$ws.Range("J51").Value = 9110.6
$ws.Range("K51").Value = 37500
$ws.Range("L51").Value = 9110.6
$ws.Range("M51").Value = -37016
$ws.Range("N51").Value = -10078.6

$ws.Range("H61").Value = 1705.375
$ws.Range("J61").Value = 1484.3334
$ws.Range("L61").Value = 4453.0002
$ws.Range("N61").Value = -4797.0002

$ws.Range("H80").Value = 723.1667
$ws.Range("I80").Value = 357.2
$ws.Range("J80").Value = 984.5714
$ws.Range("K80").Value = 1071.6
$ws.Range("L80").Value = 2953.7142
$ws.Range("M80").Value = -73.59999999999991
$ws.Range("N80").Value = -4949.7142

$ws.Range("H82").Value = 13760.75
$ws.Range("I82").Value = 3033
$ws.Range("K82").Value = 9099
$ws.Range("M82").Value = -8693

$ws.Range("H83").Value = 723.1667
$ws.Range("I83").Value = 357.2
$ws.Range("J83").Value = 984.5714
$ws.Range("K83").Value = 3214.8
$ws.Range("L83").Value = 8861.142600000001
$ws.Range("M83").Value = 1777.2
$ws.Range("N83").Value = -18845.1426

$ws.Range("H85").Value = 13760.75
$ws.Range("I85").Value = 3033
$ws.Range("K85").Value = 9099
$ws.Range("M85").Value = -7695

$ws.Range("H92").Value = 2283.2
$ws.Range("I92").Value = 2055.25
$ws.Range("J92").Value = 3195
$ws.Range("K92").Value = 2055.25
$ws.Range("L92").Value = 3195
$ws.Range("M92").Value = -807.25
$ws.Range("N92").Value = -5691

$ws.Range("H96").Value = 913.3077
$ws.Range("I96").Value = 775.55554
$ws.Range("K96").Value = 2326.66662
$ws.Range("M96").Value = -953.66662

$ws.Range("H111").Value = 2186.4167
$ws.Range("I111").Value = 1748.5
$ws.Range("J111").Value = 3062.25
$ws.Range("K111").Value = 5245.5
$ws.Range("L111").Value = 9186.75
$ws.Range("M111").Value = -2178.5
$ws.Range("N111").Value = -15320.75

$ws.Range("H123").Value = 99499.5
$ws.Range("J123").Value = 99499.5
$ws.Range("L123").Value = 99499.5
$ws.Range("N123").Value = -109299.5

$ws.Range("H132").Value = 3513.7368
$ws.Range("I132").Value = 3342.7222
$ws.Range("K132").Value = 10028.1666
$ws.Range("M132").Value = -7498.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2144.5225
$ws.Range("I32").Value = 1044.459
$ws.Range("J32").Value = 13328.5
$ws.Range("K32").Value = 1044.459
$ws.Range("L32").Value = 13328.5
$ws.Range("M32").Value = -757.4590000000001
$ws.Range("N32").Value = -13902.5

$ws.Range("H74").Value = 2436.1072
$ws.Range("I74").Value = 2441.0715
$ws.Range("K74").Value = 2441.0715
$ws.Range("M74").Value = -1567.0715

$ws.Range("H77").Value = 2436.1072
$ws.Range("I77").Value = 2441.0715
$ws.Range("K77").Value = 12205.3575
$ws.Range("M77").Value = -7837.3575

$ws.Range("H122").Value = 2413.4517
$ws.Range("I122").Value = 1600.1666
$ws.Range("K122").Value = 4800.4998
$ws.Range("M122").Value = -2350.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5519.7144
$ws.Range("I86").Value = 2254.8845
$ws.Range("J86").Value = 14951.444
$ws.Range("K86").Value = 2254.8845
$ws.Range("L86").Value = 14951.444
$ws.Range("M86").Value = -1131.8845
$ws.Range("N86").Value = -17197.444

$ws.Range("H89").Value = 5519.7144
$ws.Range("I89").Value = 2254.8845
$ws.Range("J89").Value = 14951.444
$ws.Range("K89").Value = 11274.4225
$ws.Range("L89").Value = 74757.22
$ws.Range("M89").Value = -5658.422500000001
$ws.Range("N89").Value = -85989.22

$ws.Range("H105").Value = 4286.375
$ws.Range("I105").Value = 4215.1665
$ws.Range("K105").Value = 4215.1665
$ws.Range("M105").Value = -2468.1665

$ws.Range("H107").Value = 2610.4
$ws.Range("J107").Value = 4356.5713
$ws.Range("L107").Value = 4356.5713
$ws.Range("N107").Value = -8196.5713

$ws.Range("H134").Value = 5959.1836
$ws.Range("I134").Value = 4343.657
$ws.Range("K134").Value = 13030.971
$ws.Range("M134").Value = -10495.971

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1724
$ws.Range("I99").Value = 1586.9375
$ws.Range("K99").Value = 1586.9375
$ws.Range("M99").Value = -88.9375

$ws.Range("H105").Value = 775.8889
$ws.Range("I105").Value = 654.8570999999999
$ws.Range("J105").Value = 1199.5
$ws.Range("K105").Value = 654.8570999999999
$ws.Range("L105").Value = 1199.5
$ws.Range("M105").Value = 1092.1429
$ws.Range("N105").Value = -4693.5

$ws.Range("H111").Value = 79998.5
$ws.Range("J111").Value = 79998.5
$ws.Range("L111").Value = 79998.5
$ws.Range("N111").Value = -88178.5

$ws.Range("H125").Value = 29999
$ws.Range("J125").Value = 29999
$ws.Range("L125").Value = 29999
$ws.Range("N125").Value = -34919

$ws.Range("H126").Value = 1724
$ws.Range("I126").Value = 1586.9375
$ws.Range("K126").Value = 4760.8125
$ws.Range("M126").Value = -2290.8125

$ws.Range("H134").Value = 1215.9714
$ws.Range("I134").Value = 1143.8788
$ws.Range("K134").Value = 3431.6364
$ws.Range("M134").Value = -896.6363999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 615.2
$ws.Range("I12").Value = 999
$ws.Range("J12").Value = 519.25
$ws.Range("K12").Value = 2997
$ws.Range("L12").Value = 1557.75
$ws.Range("M12").Value = -2824
$ws.Range("N12").Value = -1903.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9574.833000000001
$ws.Range("I102").Value = 13058.454
$ws.Range("J102").Value = 4100.5713
$ws.Range("K102").Value = 13058.454
$ws.Range("L102").Value = 4100.5713
$ws.Range("M102").Value = -11436.454
$ws.Range("N102").Value = -7344.5713

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 3705.818
$ws.Range("I122").Value = 3286.8125
$ws.Range("K122").Value = 9860.4375
$ws.Range("M122").Value = -7410.4375

$ws.Range("H132").Value = 3020.2
$ws.Range("I132").Value = 2745.3845
$ws.Range("J132").Value = 3530.5715
$ws.Range("K132").Value = 8236.1535
$ws.Range("L132").Value = 10591.7145
$ws.Range("M132").Value = -5706.1535
$ws.Range("N132").Value = -15651.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 443
$ws.Range("I16").Value = 454.7857
$ws.Range("K16").Value = 454.7857
$ws.Range("M16").Value = -284.7857

$ws.Range("H82").Value = 6208.5
$ws.Range("I82").Value = 1693.75
$ws.Range("K82").Value = 1693.75
$ws.Range("M82").Value = -1332.75

$ws.Range("H85").Value = 6208.5
$ws.Range("I85").Value = 1693.75
$ws.Range("K85").Value = 1693.75
$ws.Range("M85").Value = -445.75

$ws.Range("H122").Value = 5666.125
$ws.Range("I122").Value = 4168.7
$ws.Range("K122").Value = 12506.1
$ws.Range("M122").Value = -10056.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1191.125
$ws.Range("I107").Value = 1191.125
$ws.Range("K107").Value = 3573.375
$ws.Range("M107").Value = -1653.375

$ws.Range("H122").Value = 4905.45
$ws.Range("I122").Value = 3483.7058
$ws.Range("J122").Value = 12962
$ws.Range("K122").Value = 10451.1174
$ws.Range("L122").Value = 38886
$ws.Range("M122").Value = -8001.117400000001
$ws.Range("N122").Value = -43786

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
